$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date column) from row 2 to row 106:
# set the value to 45179 (2023-09-10) wherever it currently holds 45178.
for ($row = 2; $row -le 106; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value = 45179
    }
}
